$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.014.68"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.90"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.23"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6533"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.46"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2933"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07336"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.833.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6656"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.074"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008647"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "28.936.53"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.084.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "223.79"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.04%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.099"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.000"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.498"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1379"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.88"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.502"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.103"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.201"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.012"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05349"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7430"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.835"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.153"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.642"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.300.62"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.743"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.372"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8953"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.06"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.983.28"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.736"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.34%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000119"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.32%  "
$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07383"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -9.80%  "
